# Quarterly indexing esoteric bug-fix operation
#
# For each data row (a quarterly error series), a newly computed error
# value needs to be inserted at the front of the series (column B, the
# "Q0" / most-recent-quarter slot). Every existing value in that row
# slides one column to the right (Q0->Q1, Q1->Q2, ...), and the value
# that was previously in the last populated column of the row (the
# oldest quarter on record for that row) drops off the end of the
# staircase-shaped table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert into column B for each row (r=2..16)
$newVals = @{
    2  = -1.025188112727922
    3  = 0.08364543516793629
    4  = -0.1538585523806955
    5  = 0.7495351060200912
    6  = 0.03849281619118239
    7  = -0.2590580299438133
    8  = 0.01855976243503714
    9  = 0.1467044301255134
    10 = -0.1819613811903656
    11 = 0.4718454808444464
    12 = -0.08594117411414147
    13 = -0.07695400962807622
    14 = -0.5068991247689255
    15 = 0.6215838649243215
    16 = -0.2766911554241067
}

foreach ($r in 2..16) {
    # Shift existing values one column to the right, working from the
    # rightmost column back to column B so values are not clobbered
    # before they are read. Column K (11) is the last column in the
    # table, so anything already sitting there is pushed off the edge
    # and intentionally dropped.
    for ($c = 11; $c -ge 3; $c--) {
        $srcCell = $ws.Cells.Item($r, $c - 1)
        $srcVal = $srcCell.Value2
        if ($null -ne $srcVal) {
            $ws.Cells.Item($r, $c).Value2 = $srcVal
        }
    }
    # Insert the newly computed value into column B (the front of the series)
    $ws.Cells.Item($r, 2).Value2 = $newVals[$r]
}
